# daily auto push: 2025-10-07 07:26 UTC
#
# The tracker sheet (Sheet1) gets one new log row appended each day.
# Before this run the used range was A1:D73 (header row 1 + 72 data rows);
# after it, a new row 74 is appended with the latest sample:
#   A74 = "2025/10/07"  (date, stored as text - matches the existing column)
#   B74 = "火"           (weekday, text)
#   C74 = 16              (hour, number)
#   D74 = 10              (ranking, number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 74

# Column A holds date-like strings ("2025/10/07") that must stay plain text,
# exactly like every row above it. Assigning that literal straight into
# Value would get auto-recognised as a real date (Excel's usual "smart"
# typed-value coercion) and serialized as a numeric date serial instead of
# text. To avoid that without leaving a stray number-format style behind on
# the target cell, build the text in an unused scratch cell via a formula
# (which always yields a plain string result), then copy/paste just the
# value over - the destination cell keeps the sheet's default style.
$scratch = $ws.Cells.Item(1000, 26)
$scratch.Formula = '=TEXT(DATE(2025,10,7),"yyyy/mm/dd")'
$scratch.Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item($newRow, 2).Value = "火"
$ws.Cells.Item($newRow, 3).Value = 16
$ws.Cells.Item($newRow, 4).Value = 10
